$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Fecha 2021-06-21 -> now 44623, Calidad Primera, $/paquete)
$ws.Range("D2").Value = 44623
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 1800
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1900
$ws.Range("N2").Value = "$/paquete"
$ws.Range("P2").Value = 1900
$ws.Range("Q2").Value = 1

# Row 3 (now 44377, Calidad Segunda, $/docena de matas)
$ws.Range("D3").Value = 44377
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 550
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2800
$ws.Range("M3").Value = 2364
$ws.Range("N3").Value = "$/docena de matas"
$ws.Range("P3").Value = 394
$ws.Range("Q3").Value = 6

# Row 4 (now 44267, Calidad Primera)
$ws.Range("D4").Value = 44267
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = 1650
$ws.Range("P4").Value = 275

# Row 5 (now 44370, Calidad Segunda)
$ws.Range("D5").Value = 44370
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1200
$ws.Range("M5").Value = 1080
$ws.Range("P5").Value = 180
